# ---------------------------------------------------------------------------
# Updated the Code to read data from Excel Sheet
#
# Rebuilds the "Admin" sheet test-data table (new columns for login/new-user
# fields), refreshes the "TIME" sheet test-case ids (aa/bb/cc -> sequential
# Time_TestCase0N labels), and restores the view/selection state that a
# fresh edit+save pass in Excel leaves behind (active sheet back on Admin,
# per-sheet cell selections).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$admin = $wb.Worksheets.Item("Admin")
$pim   = $wb.Worksheets.Item("PIM")
$time  = $wb.Worksheets.Item("TIME")

# A stable, never-touched cell that already carries the "bordered data cell"
# look (thin box border, no fill) so new cells can inherit the exact same
# style via copy/paste-format instead of re-deriving a lookalike style.
$borderedFormatSource = $pim.Range("B1")

# ---------------------------------------------------------------------------
# Admin sheet: replace the old 5-column (A:E) mini test-case table with the
# new 10-column (A:J) table that also exercises the "add new user" fields.
# ---------------------------------------------------------------------------

# Clear the previous test-data block before laying out the new one.
$admin.Range("A1:J4").Clear() | Out-Null

# Row 1: headers, bookended by the test-case id in A1/J1.
$admin.Range("A1").Value = "Admin_TC01"
$admin.Range("B1").Value = "loginUser"
$admin.Range("C1").Value = "loginPassWord"
$admin.Range("D1").Value = "empName"
$admin.Range("E1").Value = "userRole"
$admin.Range("F1").Value = "userName"
$admin.Range("G1").Value = "newUserPassWord"
$admin.Range("H1").Value = "newUserConfirmPassWord"
$admin.Range("I1").Value = "status"
$admin.Range("J1").Value = "Admin_TC01"

# Row 2: values for test case 1.
$admin.Range("B2").Value = "Admin"
$admin.Range("C2").Value = "admin123"
$admin.Range("D2").Value = "Lokesh"
$admin.Range("E2").Value = "Admin"
$admin.Range("F2").Value = "zakir_qa"
$admin.Range("G2").Value = "rules123"
$admin.Range("H2").Value = "rules123"
$admin.Range("I2").Value = "Enabled"

# Row 3: headers for the second test case, bookended by the id in A3/E3.
$admin.Range("A3").Value = "Admin_TC02"
$admin.Range("B3").Value = "loginUser"
$admin.Range("C3").Value = "loginPassword"
$admin.Range("D3").Value = "userName"
$admin.Range("E3").Value = "Admin_TC02"

# Row 4: values for test case 2.
$admin.Range("B4").Value = "Admin"
$admin.Range("C4").Value = "admin123"
$admin.Range("D4").Value = "Lokesh"

# Formatting: every populated "data" cell gets the thin-border look (copied
# from an existing cell so the run reuses that exact style).
$borderedFormatSource.Copy() | Out-Null
$admin.Range("B1:I2").PasteSpecial(-4122) | Out-Null
$admin.Range("B3:D3").PasteSpecial(-4122) | Out-Null
$admin.Range("B4:D4").PasteSpecial(-4122) | Out-Null

# The test-case id cells (first/last column of each block) get the yellow
# highlight fill without a border.
$admin.Range("A1").Borders.LineStyle = -4142
$admin.Range("A1").Interior.Color = 65535
$idFormatSource = $admin.Range("A1")
$idFormatSource.Copy() | Out-Null
$admin.Range("J1").PasteSpecial(-4122) | Out-Null
$admin.Range("A3").PasteSpecial(-4122) | Out-Null
$admin.Range("E3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# TIME sheet: renumber the test-case id cells sequentially.
# ---------------------------------------------------------------------------

$time.Range("A1").Value = "Time_TestCase01"
$time.Range("E1").Value = "Time_TestCase01"

$time.Range("A3").Value = "Time_TestCase02"
$time.Range("E3").Value = "Time_TestCase02"

$time.Range("A5").Value = "Time_TestCase03"
$time.Range("G5").Value = "Time_TestCase03"

$time.Range("A7").Value = "Time_TestCase04"
$time.Range("E7").Value = "Time_TestCase04"

# The repeated-id cells at the end of each TIME row lose their right border
# (matching the "open" yellow box used elsewhere on the sheet).
foreach ($addr in @("E1", "E3", "G5", "E7")) {
    $time.Range($addr).Borders.Item(10).LineStyle = -4142
}

# ---------------------------------------------------------------------------
# View state: active sheet moves back to Admin, each sheet keeps its own
# remembered selection.
# ---------------------------------------------------------------------------

$pim.Range("C26").Select() | Out-Null
$time.Range("M18").Select() | Out-Null

$admin.Activate() | Out-Null
$admin.Range("G12").Select() | Out-Null
